$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 6500
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 6500
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -7468

$ws.Range("H62").Value = 1713.125
$ws.Range("I62").Value = 1713.125
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1713.125
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1089.125
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 1713.125
$ws.Range("I65").Value = 1713.125
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 8565.625
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -5445.625
$ws.Range("N65").ClearContents()

$ws.Range("H106").Value = 1845.6666
$ws.Range("I106").Value = 1455.8889
$ws.Range("K106").Value = 1455.8889
$ws.Range("M106").Value = -824.8888999999999

$ws.Range("H129").Value = 1307.6267
$ws.Range("J129").Value = 1367.5
$ws.Range("L129").Value = 4102.5
$ws.Range("N129").Value = -14102.5

$ws.Range("H137").Value = 758997.25
$ws.Range("I137").Value = 1987616.5
$ws.Range("J137").Value = 2923.8718
$ws.Range("K137").Value = 5962849.5
$ws.Range("L137").Value = 8771.615399999999
$ws.Range("M137").Value = -5960299.5
$ws.Range("N137").Value = -13871.6154

$ws.Range("H139").Value = 48296
$ws.Range("J139").Value = 48296
$ws.Range("L139").Value = 48296
$ws.Range("N139").Value = -58576

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1990.4762
$ws.Range("I61").Value = 1989.4736
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1989.4736
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1777.4736
$ws.Range("N61").Value = -2424

$ws.Range("H136").Value = 1990.4762
$ws.Range("I136").Value = 1989.4736
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5968.4208
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3418.4208
$ws.Range("N136").Value = -11100

$ws.Range("H137").Value = 40585
$ws.Range("J137").Value = 40585
$ws.Range("L137").Value = 40585
$ws.Range("N137").Value = -50785

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3794.1143
$ws.Range("I134").Value = 1518.5
$ws.Range("J134").Value = 4981.391
$ws.Range("K134").Value = 4555.5
$ws.Range("L134").Value = 14944.173
$ws.Range("M134").Value = -2020.5
$ws.Range("N134").Value = -20014.173

$ws.Range("H137").Value = 52183.125
$ws.Range("J137").Value = 52183.125
$ws.Range("L137").Value = 52183.125
$ws.Range("N137").Value = -62383.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5292828.5
$ws.Range("I16").Value = 10102376
$ws.Range("K16").Value = 10102376
$ws.Range("M16").Value = -10102089

$ws.Range("H58").Value = 2777.3872
$ws.Range("I58").Value = 1441.6666
$ws.Range("J58").Value = 7357
$ws.Range("K58").Value = 1441.6666
$ws.Range("L58").Value = 7357
$ws.Range("M58").Value = -1238.6666
$ws.Range("N58").Value = -7763

$ws.Range("H99").Value = 4388.4375
$ws.Range("I99").Value = 2202
$ws.Range("J99").Value = 5700.3
$ws.Range("K99").Value = 2202
$ws.Range("L99").Value = 5700.3
$ws.Range("M99").Value = -704
$ws.Range("N99").Value = -8696.299999999999

$ws.Range("H113").Value = 5292828.5
$ws.Range("I113").Value = 10102376
$ws.Range("K113").Value = 10102376
$ws.Range("M113").Value = -10100206

$ws.Range("H122").Value = 5569
$ws.Range("I122").Value = 4666.6665
$ws.Range("J122").Value = 6471.3335
$ws.Range("K122").Value = 13999.9995
$ws.Range("L122").Value = 19414.0005
$ws.Range("M122").Value = -11549.9995
$ws.Range("N122").Value = -24314.0005

$ws.Range("H126").Value = 4388.4375
$ws.Range("I126").Value = 2202
$ws.Range("J126").Value = 5700.3
$ws.Range("K126").Value = 6606
$ws.Range("L126").Value = 17100.9
$ws.Range("M126").Value = -4136
$ws.Range("N126").Value = -22040.9

$ws.Range("H136").Value = 2777.3872
$ws.Range("I136").Value = 1441.6666
$ws.Range("J136").Value = 7357
$ws.Range("K136").Value = 4324.9998
$ws.Range("L136").Value = 22071
$ws.Range("M136").Value = -1774.9998
$ws.Range("N136").Value = -27171

$ws.Range("H137").Value = 44465.715
$ws.Range("J137").Value = 44465.715
$ws.Range("L137").Value = 44465.715
$ws.Range("N137").Value = -54665.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1029064.9
$ws.Range("I5").Value = 3300
$ws.Range("J5").Value = 1215567.5
$ws.Range("K5").Value = 9900
$ws.Range("L5").Value = 3646702.5
$ws.Range("M5").Value = -9788
$ws.Range("N5").Value = -3646926.5

$ws.Range("H68").Value = 1236.8853
$ws.Range("J68").Value = 1474.4828
$ws.Range("L68").Value = 4423.4484
$ws.Range("N68").Value = -6045.4484

$ws.Range("H71").Value = 1236.8853
$ws.Range("J71").Value = 1474.4828
$ws.Range("L71").Value = 13270.3452
$ws.Range("N71").Value = -21382.3452

$ws.Range("H113").Value = 1761108.4
$ws.Range("I113").Value = 546.64
$ws.Range("J113").Value = 5952922
$ws.Range("K113").Value = 1639.92
$ws.Range("L113").Value = 17858766
$ws.Range("M113").Value = 530.0799999999999
$ws.Range("N113").Value = -17863106

$ws.Range("H131").Value = 785.6799999999999
$ws.Range("J131").Value = 838.093
$ws.Range("L131").Value = 2514.279
$ws.Range("N131").Value = -12594.279

$ws.Range("H135").Value = 1029064.9
$ws.Range("I135").Value = 3300
$ws.Range("J135").Value = 1215567.5
$ws.Range("K135").Value = 29700
$ws.Range("L135").Value = 10940107.5
$ws.Range("M135").Value = -27165
$ws.Range("N135").Value = -10945177.5

$ws.Range("H139").Value = 1943.2667
$ws.Range("I139").Value = 914.9
$ws.Range("K139").Value = 2744.7
$ws.Range("M139").Value = 2395.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 237.23077
$ws.Range("I2").Value = 96.40000000000001
$ws.Range("K2").Value = 96.40000000000001
$ws.Range("M2").Value = 16.59999999999999

$ws.Range("H46").Value = 25487.223
$ws.Range("J46").Value = 25516.875
$ws.Range("L46").Value = 25516.875
$ws.Range("N46").Value = -25828.875

$ws.Range("H70").Value = 8071.1177
$ws.Range("I70").Value = 6800
$ws.Range("J70").Value = 9501.125
$ws.Range("K70").Value = 6800
$ws.Range("L70").Value = 9501.125
$ws.Range("M70").Value = -6530
$ws.Range("N70").Value = -10041.125

$ws.Range("H73").Value = 8071.1177
$ws.Range("I73").Value = 6800
$ws.Range("J73").Value = 9501.125
$ws.Range("K73").Value = 6800
$ws.Range("L73").Value = 9501.125
$ws.Range("M73").Value = -5864
$ws.Range("N73").Value = -11373.125

$ws.Range("H80").Value = 2588.9
$ws.Range("I80").Value = 2384.1428
$ws.Range("J80").Value = 3066.6667
$ws.Range("K80").Value = 2384.1428
$ws.Range("L80").Value = 3066.6667
$ws.Range("M80").Value = -1386.1428
$ws.Range("N80").Value = -5062.6667

$ws.Range("H83").Value = 2588.9
$ws.Range("I83").Value = 2384.1428
$ws.Range("J83").Value = 3066.6667
$ws.Range("K83").Value = 11920.714
$ws.Range("L83").Value = 15333.3335
$ws.Range("M83").Value = -6928.714
$ws.Range("N83").Value = -25317.3335

$ws.Range("H122").Value = 6666.6665
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -31900

$ws.Range("H132").Value = 4064.7646
$ws.Range("I132").Value = 2508.6667
$ws.Range("J132").Value = 7799.4
$ws.Range("K132").Value = 7526.000100000001
$ws.Range("L132").Value = 23398.2
$ws.Range("M132").Value = -4996.000100000001
$ws.Range("N132").Value = -28458.2

$ws.Range("H137").Value = 39166.668
$ws.Range("J137").Value = 48750
$ws.Range("L137").Value = 48750
$ws.Range("N137").Value = -58950

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H136").Value = 2898.2222
$ws.Range("I136").Value = 1383.3667
$ws.Range("J136").Value = 5927.933
$ws.Range("K136").Value = 4150.1001
$ws.Range("L136").Value = 17783.799
$ws.Range("M136").Value = -1600.1001
$ws.Range("N136").Value = -22883.799

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 39994
$ws.Range("J28").Value = 39994
$ws.Range("L28").Value = 39994
$ws.Range("N28").Value = -40690

$ws.Range("H122").Value = 8562.375
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 9571.286
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 28713.858
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -33613.858

$ws.Range("H136").Value = 3189.4856
$ws.Range("I136").Value = 1378.0526
$ws.Range("J136").Value = 5340.5625
$ws.Range("K136").Value = 4134.1578
$ws.Range("L136").Value = 16021.6875
$ws.Range("M136").Value = -1584.1578
$ws.Range("N136").Value = -21121.6875
